# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.520.31'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.422.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.11'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.30'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +5.70%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.69%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.20'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.127'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0800'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.97'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.93'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.801.46'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.440.82'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.82%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.411.91'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.45'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.10%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.38'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0924'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.76'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.74%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.18%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.43%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.19'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.86%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.58'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.35'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.13%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.58%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +17.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.51'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +12.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.18'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.35%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0773'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.05%  '

$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.46%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '123.04'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.19%  '

$ws.Range("B41").Value = 'WEMIXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.23'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.78%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.109'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.12'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.70%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.948.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.77%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.11%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.50'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.68%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +9.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.79'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.07'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.08%  '

